# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# Mirrors the author's fix: the scraper was only pulling team statistics,
# not the season win/loss/tie record, so three new columns are appended
# after the existing data (through column AC) and populated with the
# team's 2011 record (102 wins, 60 losses, 0 ties) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the rest of the header row (bold, bordered,
# centered) by copying the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is the same for every player on the roster.
$wins = 102
$losses = 60
$ties = 0

$firstRow = 2
$lastRow = 44

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows $firstRow-$lastRow"
